$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title text in B1: drop the period after "6.4.2.1" in the Russian title.
# Setting a new literal text forces the shared-string table to gain a new entry
# at the end and (since the old string becomes unused) drop the old one, which
# re-indexes every other shared string reference accordingly - matching the diff.
$ws.Range("B1").Value = "6.4.2.1 Общий объем забора пресной воды "

# Update data values for year 2022 (column L)
$ws.Range("L5").Value = 8741.9
# L7 previously held a formula (=L5-L8); replace it with a plain static value
$ws.Range("L7").Value = 8483.5
$ws.Range("L14").Value = 1327.6
$ws.Range("L18").Value = 54

# Update the selected cell/range shown when the workbook is opened
$ws.Range("O2").Select()
